# Insert a new weekly price record as row 29 in the "Espinaca" sheet.
# This shifts the existing rows 29-89 down to 30-90 (Excel preserves all
# their values/formatting), and the new row 29 is populated with the
# latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 29..89 down to 30..90, opening up a blank row 29.
$ws.Rows.Item(29).Insert()

# The constant (per-series) columns are identical on every data row, so
# copy them straight from the row that used to be 29 (now 30).
$ws.Range("A29").Value2 = $ws.Range("A30").Value2
$ws.Range("B29").Value2 = $ws.Range("B30").Value2
$ws.Range("C29").Value2 = $ws.Range("C30").Value2
$ws.Range("E29").Value2 = $ws.Range("E30").Value2
$ws.Range("F29").Value2 = $ws.Range("F30").Value2
$ws.Range("G29").Value2 = $ws.Range("G30").Value2
$ws.Range("H29").Value2 = $ws.Range("H30").Value2
$ws.Range("I29").Value2 = $ws.Range("I30").Value2
$ws.Range("N29").Value2 = $ws.Range("N30").Value2
$ws.Range("O29").Value2 = $ws.Range("O30").Value2
$ws.Range("Q29").Value2 = $ws.Range("Q30").Value2
$ws.Range("R29").Value2 = $ws.Range("R30").Value2

# New observation-specific values.
$ws.Range("D29").Value2 = 44797
$ws.Range("J29").Value2 = 50
$ws.Range("K29").Value2 = 6500
$ws.Range("L29").Value2 = 7000
$ws.Range("M29").Value2 = 6700
$ws.Range("P29").Value2 = 670
